# Refresh the crypto symbol table: latest price/volume/hour snapshot.
# (GitHub Actions scheduled re-scrape; some rows also show the
# "HotbitToken" entry promoted ahead of its old slot, shifting the
# rows below it down by one.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Assigning a numeric-looking string straight to .Value lets Excel
    # infer a Number/Percentage type, which both mangles the literal
    # text (e.g. trailing zeros, "% suffix") and loses precision.
    # A leading apostrophe forces Text, and resetting the Style back to
    # Normal afterwards clears the "quote prefix" marker Excel adds so
    # the cell ends up with no explicit style, same as the source file.
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '330.92'
Set-TextValue "E2" '0.34%'
Set-TextValue "G2" '9'

# Row 3
Set-TextValue "D3" '45.48'
Set-TextValue "E3" '3.28%'
Set-TextValue "G3" '9'

# Row 4
Set-TextValue "D4" '5.612'
Set-TextValue "E4" '2.09%'
Set-TextValue "G4" '9'

# Row 5
Set-TextValue "D5" '0.08351'
Set-TextValue "E5" '4.32%'
Set-TextValue "G5" '9'

# Row 6
Set-TextValue "D6" '2.102'
Set-TextValue "E6" '6.28%'
Set-TextValue "G6" '9'

# Row 7
Set-TextValue "D7" '0.9617'
Set-TextValue "E7" '0.84%'
Set-TextValue "G7" '9'

# Row 8
Set-TextValue "D8" '2.533'
Set-TextValue "E8" '-1.59%'
Set-TextValue "G8" '9'

# Row 9
Set-TextValue "D9" '0.1157'
Set-TextValue "E9" '3.37%'
Set-TextValue "G9" '9'

# Row 10
Set-TextValue "D10" '0.1917'
Set-TextValue "E10" '2.00%'
Set-TextValue "G10" '9'

# Row 11
Set-TextValue "D11" '10.40'
Set-TextValue "E11" '-1.32%'
Set-TextValue "G11" '9'

# Row 12
Set-TextValue "D12" '0.09878'
Set-TextValue "E12" '-1.03%'
Set-TextValue "G12" '9'

# Row 13
Set-TextValue "D13" '0.04629'
Set-TextValue "E13" '-1.86%'
Set-TextValue "G13" '9'

# Row 14
Set-TextValue "D14" '0.1059'
Set-TextValue "E14" '-0.48%'
Set-TextValue "G14" '9'

# Row 15
Set-TextValue "D15" '0.001294'
Set-TextValue "E15" '1.88%'
Set-TextValue "G15" '9'

# Row 16
Set-TextValue "D16" '0.006084'
Set-TextValue "E16" '1.81%'
Set-TextValue "G16" '9'

# Row 17
$ws.Range("B17").Value = 'HotbitToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue "D17" '0.004571'
Set-TextValue "E17" '4.89%'
Set-TextValue "G17" '9'

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D18" '3.375'
Set-TextValue "E18" '0.18%'
Set-TextValue "G18" '9'

# Row 19
$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D19" '4.442'
Set-TextValue "E19" '1.44%'
Set-TextValue "G19" '9'

# Row 20
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue "D20" '0.3342'
Set-TextValue "E20" '-3.66%'
Set-TextValue "G20" '9'

# Row 21
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue "D21" '0.1394'
Set-TextValue "E21" '-2.17%'
Set-TextValue "G21" '9'

# Row 22
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue "D22" '0.2654'
Set-TextValue "E22" '2.63%'
Set-TextValue "G22" '9'

# Row 23
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue "D23" '0.04188'
Set-TextValue "E23" '2.52%'
Set-TextValue "G23" '9'

# Row 24
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue "D24" '0.001312'
Set-TextValue "E24" '3.09%'
Set-TextValue "G24" '9'

# Row 25
Set-TextValue "D25" '0.0001303'
Set-TextValue "E25" '8.60%'
Set-TextValue "G25" '9'

# Row 26
Set-TextValue "D26" '0.0003750'
Set-TextValue "E26" '0.15%'
Set-TextValue "G26" '9'

# Row 27
Set-TextValue "G27" '9'

# Row 28
Set-TextValue "G28" '9'

# Row 29
Set-TextValue "G29" '9'

# Row 30
Set-TextValue "G30" '9'

# Row 31
Set-TextValue "G31" '9'

# Row 32
Set-TextValue "G32" '9'

# Row 33
Set-TextValue "G33" '9'

# Row 34
Set-TextValue "G34" '9'

# Row 35
Set-TextValue "G35" '9'

# Row 36
Set-TextValue "G36" '9'

# Row 37
Set-TextValue "G37" '9'

# Row 38
Set-TextValue "D38" '0.02688'
Set-TextValue "E38" '3.96%'
Set-TextValue "G38" '9'

# Row 39
Set-TextValue "E39" '1.22%'
Set-TextValue "G39" '9'

# Row 40
Set-TextValue "D40" '0.007847'
Set-TextValue "E40" '3.65%'
Set-TextValue "G40" '9'

# Row 41
Set-TextValue "D41" '0.1434'
Set-TextValue "E41" '2.46%'
Set-TextValue "G41" '9'

# Row 42
Set-TextValue "D42" '0.007251'
Set-TextValue "E42" '-2.17%'
Set-TextValue "G42" '9'

# Row 43
Set-TextValue "D43" '0.002015'
Set-TextValue "E43" '0.00%'
Set-TextValue "G43" '9'

# Row 44
Set-TextValue "D44" '0.009081'
Set-TextValue "E44" '8.95%'
Set-TextValue "G44" '9'

# Row 45
Set-TextValue "D45" '0.3537'
Set-TextValue "G45" '9'

# Row 46
Set-TextValue "D46" '0.00007111'
Set-TextValue "E46" '-0.31%'
Set-TextValue "G46" '9'

# Row 47
Set-TextValue "E47" '0.26%'
Set-TextValue "G47" '9'

# Row 48
Set-TextValue "E48" '0.33%'
Set-TextValue "G48" '9'

# Row 49
Set-TextValue "D49" '0.003524'
Set-TextValue "E49" '-1.02%'
Set-TextValue "G49" '9'

# Row 50
Set-TextValue "D50" '0.003507'
Set-TextValue "E50" '-0.65%'
Set-TextValue "G50" '9'

# Row 51
Set-TextValue "E51" '0.26%'
Set-TextValue "G51" '9'
